# Corrected stupid error in spready inputs:
# C18 was a hard-coded 0.5 constant; it should be derived as the
# remainder of the 60% pours (1 - C17), matching the 4:6 method's logic.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("4-6 Pour Over Recipe")

$ws.Range("C18").Formula = "=1-C17"

# Move the active selection to C17 (matches the cell the author was
# working from/inspecting when making the fix).
$ws.Range("C17").Select()
